$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, shifting existing rows 27..41 down to 28..42.
$ws.Rows.Item(27).Insert()

# Fill the newly inserted row 27 with a new weekly data entry (same dimension values
# as the row that used to be at 27, but with a new date and a new volume).
$ws.Range("A27").Value = 7
$ws.Range("B27").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C27").Value = "Ñuble"
$ws.Range("D27").Value = 44895
$ws.Range("E27").Value = 16
$ws.Range("F27").Value = 300000000
$ws.Range("G27").Value = "Espárragos"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 1200
$ws.Range("K27").Value = 900
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = 950
$ws.Range("N27").Value = "$/kilo"
$ws.Range("O27").Value = "Región de Ñuble"
$ws.Range("P27").Value = 950
$ws.Range("Q27").Value = 1
$ws.Range("R27").Value = "Hortaliza"
